$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the claim number in B3 (shared string value changes from 1120170200907 to 1120170200908)
# Leading apostrophe forces text entry so it stays a shared string (not a number).
$ws.Range("B3").Value = "'1120170200908"

# Update the Importe value in C3 from 100 to 135
$ws.Range("C3").Value = 135

# Move the selection/cursor to C4 (was G7)
$ws.Range("C4").Select()
